$wb = $excel.ActiveWorkbook

# --- 1. Update HTSE sheet view: drop tabSelected, change selection to A1:C2 ---
$htse = $wb.Worksheets.Item("HTSE")
$htse.Range("A1:C2").Select()

# --- 2. Insert new FT sheet right after HTSE ---
$wsAfter = $wb.Worksheets.Item("HTSE")
$ft = $wb.Worksheets.Add($null, $wsAfter)
$ft.Name = "FT"

# --- 3. Column widths ---
$ft.Columns.Item(1).ColumnWidth = 27
$ft.Columns.Item(2).ColumnWidth = 34
$ft.Columns.Item(3).ColumnWidth = 11
$ft.Columns.Item(4).ColumnWidth = 14.6666666666667

# --- 4. Row heights ---
$ft.Rows.Item(1).RowHeight = 44
$ft.Rows.Item(3).RowHeight = 29

# --- 5. Cell values / formulas ---
$ft.Range("A1").Value = 'Source'
$ft.Range("B1").Value = 'Performance and cost analysis of liquid fuel production form H2 and CO2 based on the FT process'
$ft.Range("A2").Value = 'QoI'
$ft.Range("B2").Value = 'Description'
$ft.Range("C2").Value = 'Value'
$ft.Range("D2").Value = 'Unit'
$ft.Range("A3").Value = 'TDCC'
$ft.Range("B3").Value = 'Total Direct Capital Costs, sum of equipment installed costs'
$ft.Range("C3").Value = 257800644
$ft.Range("D3").Value = '$USD (2016)'
$ft.Range("A4").Value = 'Depreciable Capital Costs'
$ft.Range("B4").Value = 'Site preparation'
$ft.Range("C4").Formula = "=0.02*C3"
$ft.Range("D4").Value = '$USD (2016)'
$ft.Range("B5").Value = 'Eng and design'
$ft.Range("C5").Formula = "=0.1*C3"
$ft.Range("D5").Value = '$USD (2016)'
$ft.Range("B6").Value = 'Project contingency'
$ft.Range("C6").Formula = "=0.15*C3"
$ft.Range("D6").Value = '$USD (2016)'
$ft.Range("B7").Value = 'Catalyst first fill fee'
$ft.Range("C7").Value = 12251143
$ft.Range("D7").Value = '$USD (2016)'
$ft.Range("B8").Value = 'Upfront permitting costs'
$ft.Range("C8").Formula = "=0.15*C3"
$ft.Range("D8").Value = '$USD (2016)'
$ft.Range("A9").Value = 'Total depreciable capital costs'
$ft.Range("C9").Formula = "=SUM(C4:C8)"
$ft.Range("D9").Value = '$USD (2016)'
$ft.Range("A10").Value = 'Non-depreciable Capital Costs'
$ft.Range("B10").Value = 'Land'
$ft.Range("C10").Formula = "=10*55036"
$ft.Range("D10").Value = '$USD (2016)'
$ft.Range("A11").Value = 'TCI'
$ft.Range("B11").Value = 'Total Capital Investment'
$ft.Range("C11").Formula = "=SUM(C3,C9,C10)"
$ft.Range("D11").Value = '$USD (2016)'
$ft.Range("A12").Value = 'LC'
$ft.Range("B12").Value = 'Labor Cost'
$ft.Range("C12").Value = 9607972
$ft.Range("D12").Value = '$USD (2016)/year'
$ft.Range("A13").Value = 'Gen and admin'
$ft.Range("B13").Value = '20% LC'
$ft.Range("C13").Formula = "=0.2*C12"
$ft.Range("D13").Value = '$USD (2016)/year'
$ft.Range("A14").Value = 'Property taxes and insurance'
$ft.Range("B14").Value = '2% TCI'
$ft.Range("C14").Formula = "=0.02*C11"
$ft.Range("D14").Value = '$USD (2016)/year'
$ft.Range("A15").Value = 'Materials costs for maintenance'
$ft.Range("C15").Value = 1049006
$ft.Range("D15").Value = '$USD (2016)/year'
$ft.Range("A16").Value = 'Total Fixed Operating Costs'
$ft.Range("C16").Formula = "=SUM(C12:C15)"
$ft.Range("D16").Value = '$USD (2016)/year'
$ft.Range("A17").Value = 'None Energy material and utilities costs'
$ft.Range("C17").Value = 7085933
$ft.Range("D17").Value = '$USD (2016)/year'
$ft.Range("A18").Value = 'Total Variable Operating Costs (excl. feedstock and elec)'
$ft.Range("C18").Formula = "=C17"
$ft.Range("D18").Value = '$USD (2016)/year'

# --- 6. Styling: number format (text) + wrap text, per cell ---
$ft.Range("B1").WrapText = $true
$ft.Range("B3").NumberFormat = "@"
$ft.Range("B3").WrapText = $true
$ft.Range("B4").NumberFormat = "@"
$ft.Range("B4").WrapText = $true
$ft.Range("B5").NumberFormat = "@"
$ft.Range("B5").WrapText = $true
$ft.Range("B6").NumberFormat = "@"
$ft.Range("B6").WrapText = $true
$ft.Range("B7").NumberFormat = "@"
$ft.Range("B7").WrapText = $true
$ft.Range("B8").NumberFormat = "@"
$ft.Range("B8").WrapText = $true
$ft.Range("B9").NumberFormat = "@"
$ft.Range("B9").WrapText = $true
$ft.Range("B10").NumberFormat = "@"
$ft.Range("B10").WrapText = $true
$ft.Range("B12").NumberFormat = "@"
$ft.Range("B12").WrapText = $true
$ft.Range("B13").NumberFormat = "@"
$ft.Range("B13").WrapText = $true
$ft.Range("B14").NumberFormat = "@"
$ft.Range("B14").WrapText = $true

# --- 7. Styling: Good (green) fill for total cells ---
$ft.Range("C11").Style = "Good"
$ft.Range("C16").Style = "Good"
$ft.Range("C18").Style = "Good"

# --- 8. Styling: borders, per cell ---
$ft.Range("A2").Borders.Item(7).LineStyle = 1
$ft.Range("A2").Borders.Item(7).Weight = -4138
$ft.Range("A2").Borders.Item(8).LineStyle = 1
$ft.Range("A2").Borders.Item(8).Weight = -4138
$ft.Range("A2").Borders.Item(9).LineStyle = 1
$ft.Range("A2").Borders.Item(9).Weight = 2
$ft.Range("B2").Borders.Item(8).LineStyle = 1
$ft.Range("B2").Borders.Item(8).Weight = -4138
$ft.Range("B2").Borders.Item(9).LineStyle = 1
$ft.Range("B2").Borders.Item(9).Weight = 2
$ft.Range("C2").Borders.Item(8).LineStyle = 1
$ft.Range("C2").Borders.Item(8).Weight = -4138
$ft.Range("C2").Borders.Item(9).LineStyle = 1
$ft.Range("C2").Borders.Item(9).Weight = 2
$ft.Range("D2").Borders.Item(10).LineStyle = 1
$ft.Range("D2").Borders.Item(10).Weight = -4138
$ft.Range("D2").Borders.Item(8).LineStyle = 1
$ft.Range("D2").Borders.Item(8).Weight = -4138
$ft.Range("D2").Borders.Item(9).LineStyle = 1
$ft.Range("D2").Borders.Item(9).Weight = 2
$ft.Range("A3").Borders.Item(7).LineStyle = 1
$ft.Range("A3").Borders.Item(7).Weight = -4138
$ft.Range("A3").Borders.Item(8).LineStyle = 1
$ft.Range("A3").Borders.Item(8).Weight = 2
$ft.Range("B3").Borders.Item(8).LineStyle = 1
$ft.Range("B3").Borders.Item(8).Weight = 2
$ft.Range("C3").Borders.Item(8).LineStyle = 1
$ft.Range("C3").Borders.Item(8).Weight = 2
$ft.Range("D3").Borders.Item(10).LineStyle = 1
$ft.Range("D3").Borders.Item(10).Weight = -4138
$ft.Range("D3").Borders.Item(8).LineStyle = 1
$ft.Range("D3").Borders.Item(8).Weight = 2
$ft.Range("A4").Borders.Item(7).LineStyle = 1
$ft.Range("A4").Borders.Item(7).Weight = -4138
$ft.Range("D4").Borders.Item(10).LineStyle = 1
$ft.Range("D4").Borders.Item(10).Weight = -4138
$ft.Range("A5").Borders.Item(7).LineStyle = 1
$ft.Range("A5").Borders.Item(7).Weight = -4138
$ft.Range("D5").Borders.Item(10).LineStyle = 1
$ft.Range("D5").Borders.Item(10).Weight = -4138
$ft.Range("A6").Borders.Item(7).LineStyle = 1
$ft.Range("A6").Borders.Item(7).Weight = -4138
$ft.Range("D6").Borders.Item(10).LineStyle = 1
$ft.Range("D6").Borders.Item(10).Weight = -4138
$ft.Range("A7").Borders.Item(7).LineStyle = 1
$ft.Range("A7").Borders.Item(7).Weight = -4138
$ft.Range("D7").Borders.Item(10).LineStyle = 1
$ft.Range("D7").Borders.Item(10).Weight = -4138
$ft.Range("A8").Borders.Item(7).LineStyle = 1
$ft.Range("A8").Borders.Item(7).Weight = -4138
$ft.Range("D8").Borders.Item(10).LineStyle = 1
$ft.Range("D8").Borders.Item(10).Weight = -4138
$ft.Range("A9").Borders.Item(7).LineStyle = 1
$ft.Range("A9").Borders.Item(7).Weight = -4138
$ft.Range("D9").Borders.Item(10).LineStyle = 1
$ft.Range("D9").Borders.Item(10).Weight = -4138
$ft.Range("A10").Borders.Item(7).LineStyle = 1
$ft.Range("A10").Borders.Item(7).Weight = -4138
$ft.Range("D10").Borders.Item(10).LineStyle = 1
$ft.Range("D10").Borders.Item(10).Weight = -4138
$ft.Range("A11").Borders.Item(7).LineStyle = 1
$ft.Range("A11").Borders.Item(7).Weight = -4138
$ft.Range("A11").Borders.Item(9).LineStyle = 1
$ft.Range("A11").Borders.Item(9).Weight = 2
$ft.Range("B11").Borders.Item(9).LineStyle = 1
$ft.Range("B11").Borders.Item(9).Weight = 2
$ft.Range("C11").Borders.Item(9).LineStyle = 1
$ft.Range("C11").Borders.Item(9).Weight = 2
$ft.Range("D11").Borders.Item(10).LineStyle = 1
$ft.Range("D11").Borders.Item(10).Weight = -4138
$ft.Range("D11").Borders.Item(9).LineStyle = 1
$ft.Range("D11").Borders.Item(9).Weight = 2
$ft.Range("A12").Borders.Item(7).LineStyle = 1
$ft.Range("A12").Borders.Item(7).Weight = -4138
$ft.Range("A12").Borders.Item(8).LineStyle = 1
$ft.Range("A12").Borders.Item(8).Weight = 2
$ft.Range("B12").Borders.Item(8).LineStyle = 1
$ft.Range("B12").Borders.Item(8).Weight = 2
$ft.Range("C12").Borders.Item(8).LineStyle = 1
$ft.Range("C12").Borders.Item(8).Weight = 2
$ft.Range("D12").Borders.Item(10).LineStyle = 1
$ft.Range("D12").Borders.Item(10).Weight = -4138
$ft.Range("D12").Borders.Item(8).LineStyle = 1
$ft.Range("D12").Borders.Item(8).Weight = 2
$ft.Range("A13").Borders.Item(7).LineStyle = 1
$ft.Range("A13").Borders.Item(7).Weight = -4138
$ft.Range("D13").Borders.Item(10).LineStyle = 1
$ft.Range("D13").Borders.Item(10).Weight = -4138
$ft.Range("A14").Borders.Item(7).LineStyle = 1
$ft.Range("A14").Borders.Item(7).Weight = -4138
$ft.Range("D14").Borders.Item(10).LineStyle = 1
$ft.Range("D14").Borders.Item(10).Weight = -4138
$ft.Range("A15").Borders.Item(7).LineStyle = 1
$ft.Range("A15").Borders.Item(7).Weight = -4138
$ft.Range("D15").Borders.Item(10).LineStyle = 1
$ft.Range("D15").Borders.Item(10).Weight = -4138
$ft.Range("A16").Borders.Item(7).LineStyle = 1
$ft.Range("A16").Borders.Item(7).Weight = -4138
$ft.Range("D16").Borders.Item(10).LineStyle = 1
$ft.Range("D16").Borders.Item(10).Weight = -4138
$ft.Range("A17").Borders.Item(7).LineStyle = 1
$ft.Range("A17").Borders.Item(7).Weight = -4138
$ft.Range("A17").Borders.Item(8).LineStyle = 1
$ft.Range("A17").Borders.Item(8).Weight = 2
$ft.Range("B17").Borders.Item(8).LineStyle = 1
$ft.Range("B17").Borders.Item(8).Weight = 2
$ft.Range("C17").Borders.Item(8).LineStyle = 1
$ft.Range("C17").Borders.Item(8).Weight = 2
$ft.Range("D17").Borders.Item(10).LineStyle = 1
$ft.Range("D17").Borders.Item(10).Weight = -4138
$ft.Range("D17").Borders.Item(8).LineStyle = 1
$ft.Range("D17").Borders.Item(8).Weight = 2
$ft.Range("A18").Borders.Item(7).LineStyle = 1
$ft.Range("A18").Borders.Item(7).Weight = -4138
$ft.Range("A18").Borders.Item(9).LineStyle = 1
$ft.Range("A18").Borders.Item(9).Weight = -4138
$ft.Range("B18").Borders.Item(9).LineStyle = 1
$ft.Range("B18").Borders.Item(9).Weight = -4138
$ft.Range("C18").Borders.Item(9).LineStyle = 1
$ft.Range("C18").Borders.Item(9).Weight = -4138
$ft.Range("D18").Borders.Item(10).LineStyle = 1
$ft.Range("D18").Borders.Item(10).Weight = -4138
$ft.Range("D18").Borders.Item(9).LineStyle = 1
$ft.Range("D18").Borders.Item(9).Weight = -4138

# --- 9. Finalize FT sheet view: activeCell I10, select I10, activate tab ---
$ft.Range("I10").Select()
$ft.Activate()

Write-Output "done"